# "Generate Report for Handoff"
#
# A handoff xliff report was (re)generated for the
# "2257adb3-2956-4548-b452-34068782e39a" entry. This updates the
# "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Latest Handoff Datetime" for that file's row on each per-locale
# sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$fileKey = "2257adb3-2956-4548-b452-34068782e39a"

function Find-RowByFileName($ws, $col) {
    $usedRows = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $usedRows; $r++) {
        $v = $ws.Cells.Item($r, $col).Value2
        if ($v -ne $null -and $v.ToString().StartsWith($fileKey)) {
            return $r
        }
    }
    return -1
}

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$rowOverview = Find-RowByFileName $wsOverview 1
if ($rowOverview -gt 0) {
    $wsOverview.Cells.Item($rowOverview, 7).Value = "2017-02-17 06:51:17"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$rowZhCn = Find-RowByFileName $wsZhCn 1
if ($rowZhCn -gt 0) {
    $wsZhCn.Cells.Item($rowZhCn, 8).Value = "2017-02-17 06:50:59"
}

# --- de-de sheet: column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$rowDeDe = Find-RowByFileName $wsDeDe 1
if ($rowDeDe -gt 0) {
    $wsDeDe.Cells.Item($rowDeDe, 8).Value = "2017-02-17 06:51:17"
}

Write-Output "Updated rows -> Overview:$rowOverview zh-cn:$rowZhCn de-de:$rowDeDe"
